$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 'Vega Modelo de Temuco'
$ws.Range("C6").Value = 'La Araucanía'
$ws.Range("D6").Value = 45092
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100112010
$ws.Range("G6").Value = 'Achicoria'
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 125
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("N6").Value = '$/caja 18 unidades'
$ws.Range("O6").Value = 'Región Metropolitana'
$ws.Range("P6").Value = 389
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = 'Hortaliza'
